# Rename "AddressBook" (domain model class) to "VoluncheerBook" throughout
# the UML-ish diagram on slide 1, per commit "Update developer guide
# including related diagrams".
#
# Four shapes on the slide carry the old "Address" naming and need the
# substring "Address" swapped for "Voluncheer" inside their text:
#   - "Rectangle 62"  (":AddressBookParser"   -> ":VoluncheerBookParser")
#   - "TextBox 78"    ("undoAddressBook()"    -> "undoVoluncheerBook()")
#   - "Rectangle 62"  (":VersionedAddressBook"-> ":VersionedVoluncheerBook")
#   - "TextBox 87"    ("resetData(ReadOnlyAddressBook)"
#                                             -> "resetData(ReadOnlyVoluncheerBook)")
#
# Each edit below targets only the exact characters that change (via
# TextRange.Characters(start, length)) so the untouched runs/paragraphs
# keep their original formatting, and any auto-fit box just resizes
# itself the way PowerPoint would.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-SubstringText {
    # NOTE: this interpreter only binds positional parameters reliably, so
    # call this as: Set-SubstringText <shape> <oldSubstring> <newSubstring>
    param($shape, [string]$old, [string]$new)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        throw "substring not found"
    }
    $sub = $tr.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    $text = $sh.TextFrame.TextRange.Text
    if ($text -notlike "*Address*") { continue }

    if ($text -like ":Address*BookParser*") {
        # ":AddressBookParser" (rendered on two lines) - split "Address" out
        # of the single ":Address" run -> ":" + "Voluncheer"
        Set-SubstringText $sh "Address" "Voluncheer"
    }
    elseif ($text -eq "undoAddressBook()") {
        Set-SubstringText $sh "AddressBook" "VoluncheerBook"
    }
    elseif ($text -eq ":VersionedAddressBook") {
        Set-SubstringText $sh "AddressBook" "VoluncheerBook"
    }
    elseif ($text -eq "resetData(ReadOnlyAddressBook)") {
        Set-SubstringText $sh "ReadOnlyAddressBook" "ReadOnlyVoluncheerBook"
    }
}
